{"js": "// Update the date line and the 25 division-problem cells in the practice\n// table to the new day's values. Text is replaced in place (via each\n// paragraph's Range) so existing run/paragraph formatting (fonts, size,\n// alignment) is preserved.\n\n// 1) Update the heading date.\nconst dateParas = context.document.body.paragraphs;\ndateParas.load(\"items\");\nawait context.sync();\n\nconst oldDate = \"2024-08-14 Wednesday\";\nconst newDate = \"2024-08-15 Thursday\";\nfor (const para of dateParas.items) {\n  para.load(\"text\");\n}\nawait context.sync();\nfor (const para of dateParas.items) {\n  if (para.text === oldDate) {\n    para.getRange().insertText(newDate, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n\n// 2) Update the table of division problems. The table has 20 rows (5\n// \"data\" rows holding the problems, each followed by 3 blank spacer\n// rows) and 5 columns. Replace each data cell's text with the new value,\n// keyed by (row, column) so duplicate old values (e.g. \"74\u00f74=18, 2\"\n// appearing twice) are handled unambiguously.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// [rowIndex, colIndex, newText]\nconst updates = [\n  [0, 0, \"74\u00f78=9, 2\"],\n  [0, 1, \"26\u00f75=5, 1\"],\n  [0, 2, \"35\u00f72=17, 1\"],\n  [0, 3, \"98\u00f78=12, 2\"],\n  [0, 4, \"96\u00f76=16, 0\"],\n\n  [4, 0, \"69\u00f74=17, 1\"],\n  [4, 1, \"32\u00f72=16, 0\"],\n  [4, 2, \"29\u00f79=3, 2\"],\n  [4, 3, \"84\u00f77=12, 0\"],\n  [4, 4, \"50\u00f76=8, 2\"],\n\n  [8, 0, \"65\u00f72=32, 1\"],\n  [8, 1, \"57\u00f77=8, 1\"],\n  [8, 2, \"44\u00f74=11, 0\"],\n  [8, 3, \"68\u00f74=17, 0\"],\n  [8, 4, \"82\u00f78=10, 2\"],\n\n  [12, 0, \"41\u00f74=10, 1\"],\n  [12, 1, \"53\u00f76=8, 5\"],\n  [12, 2, \"70\u00f75=14, 0\"],\n  [12, 3, \"94\u00f72=47, 0\"],\n  [12, 4, \"89\u00f74=22, 1\"],\n\n  [16, 0, \"19\u00f74=4, 3\"],\n  [16, 1, \"77\u00f73=25, 2\"],\n  [16, 2, \"44\u00f72=22, 0\"],\n  [16, 3, \"23\u00f76=3, 5\"],\n  [16, 4, \"59\u00f73=19, 2\"],\n];\n\nconst cellParas = [];\nfor (const [r, c, text] of updates) {\n  const cell = table.getCell(r, c);\n  cell.body.paragraphs.load(\"items\");\n  cellParas.push({ paragraphs: cell.body.paragraphs, text });\n}\nawait context.sync();\n\nfor (const { paragraphs, text } of cellParas) {\n  const para = paragraphs.items[0];\n  para.getRange().insertText(text, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division-problem cells in the practice\n# table to the new day's values. Cell/Range.Text assignment edits the\n# run's text in place, preserving existing formatting (fonts, size,\n# paragraph alignment).\n\n$d = $word.ActiveDocument\n\n# 1) Update the heading date (first paragraph in the document).\n$d.Paragraphs.Item(1).Range.Text = \"2024-08-15 Thursday\"\n\n# 2) Update the table of division problems. The table has 20 rows (5\n# \"data\" rows holding the problems, each followed by 3 blank spacer\n# rows) and 5 columns. Cells are addressed by (row, column) -- 1-based,\n# per the Word COM object model -- so duplicate old values (e.g.\n# \"74\u00f74=18, 2\" appearing twice) are handled unambiguously.\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"74\u00f78=9, 2\"\n$t.Cell(1, 2).Range.Text = \"26\u00f75=5, 1\"\n$t.Cell(1, 3).Range.Text = \"35\u00f72=17, 1\"\n$t.Cell(1, 4).Range.Text = \"98\u00f78=12, 2\"\n$t.Cell(1, 5).Range.Text = \"96\u00f76=16, 0\"\n\n$t.Cell(5, 1).Range.Text = \"69\u00f74=17, 1\"\n$t.Cell(5, 2).Range.Text = \"32\u00f72=16, 0\"\n$t.Cell(5, 3).Range.Text = \"29\u00f79=3, 2\"\n$t.Cell(5, 4).Range.Text = \"84\u00f77=12, 0\"\n$t.Cell(5, 5).Range.Text = \"50\u00f76=8, 2\"\n\n$t.Cell(9, 1).Range.Text = \"65\u00f72=32, 1\"\n$t.Cell(9, 2).Range.Text = \"57\u00f77=8, 1\"\n$t.Cell(9, 3).Range.Text = \"44\u00f74=11, 0\"\n$t.Cell(9, 4).Range.Text = \"68\u00f74=17, 0\"\n$t.Cell(9, 5).Range.Text = \"82\u00f78=10, 2\"\n\n$t.Cell(13, 1).Range.Text = \"41\u00f74=10, 1\"\n$t.Cell(13, 2).Range.Text = \"53\u00f76=8, 5\"\n$t.Cell(13, 3).Range.Text = \"70\u00f75=14, 0\"\n$t.Cell(13, 4).Range.Text = \"94\u00f72=47, 0\"\n$t.Cell(13, 5).Range.Text = \"89\u00f74=22, 1\"\n\n$t.Cell(17, 1).Range.Text = \"19\u00f74=4, 3\"\n$t.Cell(17, 2).Range.Text = \"77\u00f73=25, 2\"\n$t.Cell(17, 3).Range.Text = \"44\u00f72=22, 0\"\n$t.Cell(17, 4).Range.Text = \"23\u00f76=3, 5\"\n$t.Cell(17, 5).Range.Text = \"59\u00f73=19, 2\"\n"}
